$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for nombre_aides (C), nombre_entreprises (D), montant_total (E)
# for the 2022-05-23 data refresh.
$updates = @(
    @{ Row = 70;  C = 15731;  D = 2768;  E = 24673307 },
    @{ Row = 79;  C = 116592; D = 22734; E = 447357508 },
    @{ Row = 91;  C = 151119; D = 24834; E = 482282585 },
    @{ Row = 92;  C = 409069; D = 70906; E = 1594718796 },
    @{ Row = 93;  C = 209545; D = 34261; E = 1308331272 },
    @{ Row = 94;  C = 94175;  D = 13795; E = 916858080 },
    @{ Row = 95;  C = 50752;  D = 6982;  E = 931855831 },
    @{ Row = 96;  C = 17256;  D = 2565;  E = 790601490 },
    @{ Row = 97;  C = 2156;   D = 375;   E = 214088295 },
    @{ Row = 104; C = 135232; D = 23286; E = 272164522 },
    @{ Row = 114; C = 3801;   D = 699;   E = 9113119 },
    @{ Row = 115; C = 11693;  D = 2248;  E = 32955791 },
    @{ Row = 165; C = 83803;  D = 17113; E = 354981515 },
    @{ Row = 167; C = 12218;  D = 2416;  E = 105743177 },
    @{ Row = 168; C = 6205;   D = 1058;  E = 100553127 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
}
